$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column F ("Address"); existing F (District) shifts to G
$ws.Columns.Item(6).Insert()

$ws.Cells.Item(2, 6).Value = "Address"

$addresses = @{
    3 = "S M V V S High SchoolAfzalpur"
    4 = "Shri Gurudev Janate TrustHigh School AdnoorAfzalpur"
    5 = "S S High SchoolHire – JewargiAfzalpur"
    6 = "Nutan High SchoolMallabadAfzalpur"
    7 = "Shri Mahanteswar High School MahantapurChinamaglraAfzalpura"
    8 = "G H S TellurAfzalpur"
    9 = "Govt. High School Gobbur (B) Afzalpur"
    10 = "G H S Deval GanagapurAfzalpur"
    11 = "Govt. Urdu High School SulepethChincholi"
    12 = "Govt High SchoolDegalmadiChincholi"
    13 = "G H S KalagiChittapur"
    14 = "G H S Bhimalli"
    15 = "Govt. High School Ferozabad"
    16 = "G H S JawalgeraSindhanur"
    17 = "Govt High SchoolSindhanur"
    18 = "Adarsha Vidyalaya Lingasugur"
    19 = "G H S YlagattaLingasugur"
    20 = "G P U C H S GuruguntaLingasugur"
    21 = "G H S TalekhanLingasugur"
    22 = "G H S MatturLingasugur"
    23 = "G B H S SirwarManvi"
    24 = "G H S (Girls) SirwarManvi"
    25 = "G H S AmeenagadManvi"
    26 = "G H S ChinchodiDeodurga"
    27 = "G H S BunkaladoddiDeodurga"
    28 = "G H S MudboolShahapur"
    29 = "G H S AnabiShahapur"
    30 = "G H S KhanapurShahapur"
    31 = "Govt. Girls Junior CollegeShahapur"
    32 = "Vidyaranaya SchoolShahapur"
    33 = "Govt. H S DarshanapurShahapur"
    34 = "Govt High School ShirwalShahapur"
    35 = "Govt. Girls P U CollegeShahapur"
    37 = "G H S HegganadoddiShorapur"
    38 = "G H S BenakanahalleShorapur"
    39 = "Govt. Junior BoysCollege Gurumatkal"
    41 = "Govt. High School Honagera"
    42 = "G G H S Gurmatkal"
    43 = "G G H S Gurmatkal"
    44 = "Govt P U CollegeGurumatkal"
    45 = "Govt. Girls High School HunasagiShorapur"
    46 = "G H S Hattikani"
    47 = "M R M H School"
}

foreach ($row in $addresses.Keys) {
    $ws.Cells.Item([int]$row, 6).Value = $addresses[$row]
}
